$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 69; this shifts the existing rows 69-114
# (and their data) down to rows 70-115, which is exactly what the target
# diff shows (every row N>=70 ends up with the data that used to live in
# row N-1, and the former row 114 becomes row 115).
$ws.Rows.Item(69).Insert()

# Populate the newly inserted row 69 with the new data point.
$ws.Range("A69").Value = 4
$ws.Range("B69").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C69").Value = "Los Lagos"
$ws.Range("D69").Value = 44873
$ws.Range("E69").Value = 10
$ws.Range("F69").Value = 100112026
$ws.Range("G69").Value = "Haba"
$ws.Range("H69").Value = "Sin especificar"
$ws.Range("I69").Value = "Primera"
$ws.Range("J69").Value = 180
$ws.Range("K69").Value = 12000
$ws.Range("L69").Value = 12000
$ws.Range("M69").Value = 12000
$ws.Range("N69").Value = "$/saco 25 kilos"
$ws.Range("O69").Value = "Región del Maule"
$ws.Range("P69").Value = 480
$ws.Range("Q69").Value = 25
$ws.Range("R69").Value = "Hortaliza"
